$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- configuration -------------------------------------------------
$maxRow = 82
$rowMap = @(0,0,2,3,4,5,8,6,7,12,9,11,10,13,15,16,14,17,18,19,20,21,22,23,24,25,26,48,49,59,67,77,38,76,30,35,33,32,46,52,54,51,81,78,41,65,42,55,57,80,61,60,68,79,37,36,63,62,72,66,64,82,29,50,70,56,69,71,27,28,53,73,74,31,47,58,75,34,39,40,43,45,44)   # index = old row number, value = new row number (0 = unused)
$uNewRows = @(2)
$zNewRows = @(2,6)
$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R")

# --- read the current ('before') data for every data row -----------
$source = @{}
for ($r = 2; $r -le $maxRow; $r++) {
    $rowData = @{}
    foreach ($c in $cols) {
        $addr = "$c$r"
        $rowData[$c] = $ws.Range($addr).Value2
    }
    $source[$r] = $rowData
}

# --- clear the old link-formula columns (S,T,U,V,W,X,Y,Z) ----------
$formulaCols = @("S","T","U","V","W","X","Y","Z")
for ($r = 2; $r -le $maxRow; $r++) {
    foreach ($c in $formulaCols) {
        $ws.Range("$c$r").ClearContents()
    }
}

# --- write the reordered data back into A:R -------------------------
for ($old = 2; $old -le $maxRow; $old++) {
    $new = $rowMap[$old]
    $rowData = $source[$old]
    foreach ($c in $cols) {
        $ws.Range("$c$new").Value2 = $rowData[$c]
    }
    # the 'last changed' column is bumped for every row
    $ws.Range("C$new").Value2 = 46064
}

# --- rebuild the hyperlink formulas for the first 15 data rows ------
$baseUrl = "https://klasma.github.io/Logging_2584"
for ($r = 2; $r -le 16; $r++) {
    $a = $ws.Range("A$r").Value2
    $ws.Range("S$r").Formula = "=HYPERLINK(""$baseUrl/artfynd/$a artfynd.xlsx"", ""$a"")"
    $ws.Range("T$r").Formula = "=HYPERLINK(""$baseUrl/kartor/$a karta.png"", ""$a"")"
    if ($uNewRows -contains $r) {
        $ws.Range("U$r").Formula = "=HYPERLINK(""$baseUrl/knärot/$a karta knärot.png"", ""$a"")"
    }
    $ws.Range("V$r").Formula = "=HYPERLINK(""$baseUrl/klagomål/$a FSC-klagomål.docx"", ""$a"")"
    $ws.Range("W$r").Formula = "=HYPERLINK(""$baseUrl/klagomålsmail/$a FSC-klagomål mail.docx"", ""$a"")"
    $ws.Range("X$r").Formula = "=HYPERLINK(""$baseUrl/tillsyn/$a tillsynsbegäran.docx"", ""$a"")"
    $ws.Range("Y$r").Formula = "=HYPERLINK(""$baseUrl/tillsynsmail/$a tillsynsbegäran mail.docx"", ""$a"")"
    if ($zNewRows -contains $r) {
        $ws.Range("Z$r").Formula = "=HYPERLINK(""$baseUrl/fåglar/$a prioriterade fågelarter.docx"", ""$a"")"
    }
}

Write-Host "Reordering complete."
